$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

# Title shape on this slide currently holds "Below" + " " + "section-level"
# split across three separate runs. Collapse them into a single run with
# the combined text "Below section-level".
$shp = $s.Shapes.Item(1)
$tr = $shp.TextFrame.TextRange

# Route the text-range update through an intermediate value that shares no
# prefix with the current text so the host rewrites the whole paragraph
# into a single run (re-using the first run's existing, empty <a:rPr/>)
# instead of only patching/appending to the trailing run.
$tr.Text = "."
$tr.Text = "Below section-level"
